# Atualização de bases das ligas, do dia: 11-04-2024 às 00:31
# Swap the data (all columns except the leading row-index column A)
# between the following row pairs: 20<->21, 33<->34, 58<->59.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($ws, $row1, $row2, $colStart, $colEnd) {
    $vals1 = @()
    $vals2 = @()
    for ($col = $colStart; $col -le $colEnd; $col++) {
        $vals1 += ,$ws.Cells.Item($row1, $col).Value2
        $vals2 += ,$ws.Cells.Item($row2, $col).Value2
    }
    $n = $colEnd - $colStart + 1
    for ($i = 0; $i -lt $n; $i++) {
        $col = $colStart + $i
        $ws.Cells.Item($row1, $col).Value2 = $vals2[$i]
        $ws.Cells.Item($row2, $col).Value2 = $vals1[$i]
    }
}

# Columns B (2) through AC (29) hold the record data; column A is just
# the sequential id and must stay untouched.
$colStart = 2
$colEnd = 29

Swap-RowData $ws 20 21 $colStart $colEnd
Swap-RowData $ws 33 34 $colStart $colEnd
Swap-RowData $ws 58 59 $colStart $colEnd
